# Scheduled-runner update: refresh market-price-derived columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the per-job Profits sheets.
$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# ---- ALC ----
Set-Row "ALC" 100 @{ H = 1200.75;     I = 1086.5714; K = 1086.5714; M = -545.5714 }
Set-Row "ALC" 114 @{ H = 42805.332;   J = 42805.332; L = 42805.332; N = -51483.332 }
Set-Row "ALC" 124 @{ H = 46780;       J = 46780;     L = 46780;     N = -56600 }
Set-Row "ALC" 128 @{ H = 51999;       J = 51999;     L = 51999;     N = -61959 }
Set-Row "ALC" 130 @{ H = 0;           J = 0;         L = 0;         N = "" }

# ---- ARM ----
Set-Row "ARM" 32  @{ H = 10630.226;   I = 10256.808; K = 10256.808; M = -9969.808000000001 }
Set-Row "ARM" 74  @{ H = 1476.3;      I = 1325.3864; K = 1325.3864; M = -451.3864000000001 }
Set-Row "ARM" 77  @{ H = 1476.3;      I = 1325.3864; K = 6626.932000000001; M = -2258.932000000001 }
Set-Row "ARM" 109 @{ H = 38764.668;   J = 38764.668; L = 38764.668; N = -41538.668 }
Set-Row "ARM" 111 @{ H = 49447;       J = 49447;     L = 49447;     N = -57627 }
Set-Row "ARM" 113 @{ H = 46746;       J = 46746;     L = 46746;     N = -55424 }
Set-Row "ARM" 121 @{ H = 37490.25;    J = 37490.25;  L = 37490.25;  N = -40984.25 }
Set-Row "ARM" 125 @{ H = 48178.25;    J = 48178.25;  L = 48178.25;  N = -58018.25 }
Set-Row "ARM" 132 @{ H = 8334786.5;   J = 2520;      L = 7560;      N = -12620 }

# ---- BSM ----
Set-Row "BSM" 111 @{ H = 47686; J = 47686; L = 47686; N = -55866 }
Set-Row "BSM" 139 @{ H = 33000; J = 33000; L = 33000; N = -43280 }

# ---- CRP ----
Set-Row "CRP" 20  @{ H = 49887.5;    J = 49887.5;    L = 49887.5;    N = -50359.5 }
Set-Row "CRP" 30  @{ H = 49887.5;    J = 49887.5;    L = 49887.5;    N = -50069.5 }
Set-Row "CRP" 110 @{ H = 46346.5;    J = 46346.5;    L = 46346.5;    N = -54526.5 }
Set-Row "CRP" 111 @{ H = 40364;      J = 40364;      L = 40364;      N = -48544 }
Set-Row "CRP" 112 @{ H = 47694;      J = 47694;      L = 47694;      N = -50648 }
Set-Row "CRP" 119 @{ H = 52637.668;  J = 52637.668;  L = 52637.668;  N = -62313.668 }
Set-Row "CRP" 128 @{ H = 49887.5;    J = 49887.5;    L = 49887.5;    N = -59847.5 }

# ---- GSM ----
Set-Row "GSM" 116 @{ H = 48742; J = 48742; L = 48742; N = -57920 }
Set-Row "GSM" 119 @{ H = 48761; J = 48761; L = 48761; N = -58437 }
Set-Row "GSM" 124 @{ H = 41780; J = 41780; L = 41780; N = -51600 }
Set-Row "GSM" 128 @{ H = 39780; J = 39780; L = 39780; N = -49740 }
Set-Row "GSM" 132 @{ H = 2898.4517; I = 1924.1177; J = 4081.5715; K = 5772.3531; L = 12244.7145; M = -3242.3531; N = -17304.7145 }

# ---- LTW ----
Set-Row "LTW" 108 @{ H = 48626;      J = 48626;      L = 48626;      N = -56306 }
Set-Row "LTW" 111 @{ H = 46253;      J = 46253;      L = 46253;      N = -54433 }
Set-Row "LTW" 112 @{ H = 43264;      J = 43264;      L = 43264;      N = -46218 }
Set-Row "LTW" 117 @{ H = 43388;      J = 43388;      L = 43388;      N = -52566 }
Set-Row "LTW" 120 @{ H = 47767.25;   J = 47767.25;   L = 47767.25;   N = -57443.25 }
Set-Row "LTW" 123 @{ H = 0;          J = 0;          L = 0;          N = "" }
Set-Row "LTW" 132 @{ H = 2494.7036;  I = 1819.58;    K = 5458.74;    M = -2928.74 }

# ---- WVR ----
Set-Row "WVR" 108 @{ H = 47626; J = 47626; L = 47626; N = -55306 }
Set-Row "WVR" 110 @{ H = 48507; J = 48507; L = 48507; N = -56687 }
Set-Row "WVR" 112 @{ H = 0;     J = 0;     L = 0;     N = "" }
Set-Row "WVR" 116 @{ H = 49680; J = 49680; L = 49680; N = -58858 }
Set-Row "WVR" 121 @{ H = 41156; J = 41156; L = 41156; N = -44650 }
Set-Row "WVR" 125 @{ H = 40684; J = 40684; L = 40684; N = -50524 }
